$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column -> new value mapping (identical update applied to rows 2 and 3)
$updates = @{
    "K"  = -2.25
    "U"  = 0.592
    "V"  = 0.02152727272727272
    "W"  = -0.1461038961038961
    "X"  = 0.09156073182779319
    "Y"  = -0.2376646279316893
    "AA" = -0.09648633575013943
    "AB" = 0.08492258372912596
    "AC" = -0.1814089194792654
    "AD" = 4.75
    "AF" = 4.75
    "AG" = 4.158
    "AH" = 0.1472868217054264
    "AI" = 0.2454780361757106
    "AJ" = 0.1313412091730368
    "AK" = 0.2216654227529588
    "AL" = 0.986
    "AM" = 0.986
    "AN" = -2.844311377245509
    "AO" = -1.754563894523327
    "AP" = -2.489820359281437
    "AQ" = -1.754563894523327
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
